# Insert a new row for an additional responsible teacher (Docentes responsáveis:)
# right below the existing "3577649 - Carlos Angelo Nunes" row, shifting all
# subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 14 (pushes old rows 14-24 down to 15-25),
# inheriting formatting from the row above as Excel normally does.
$ws.Rows(14).Insert()

# Populate the new row with the new teacher's name in both the "current" (B)
# and "modified" (C) columns, matching the pattern used by the existing
# teacher row above it.
$ws.Range("B14").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C14").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
